$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the "Save" header in H1, copying the existing header style (bold,
# centered, bordered) from an existing header cell so it matches B1:G1.
$ws.Range("H1").Value = "Save"
$ws.Range("B1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Fill in the "Save" column values for each data row.
$values = @(0, 0, 0, 0, 1, 0, 0, 1, 1, 1)
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
